$d = $word.ActiveDocument
$wdYellow = 7

# 1. "Can train only in a BoxingGym." -> highlight whole sentence (runs only,
#    not the paragraph mark).
$range = $d.Content
$found = $range.Find.Execute("Can train only in a BoxingGym.")
if ($found) {
    $range.Font.HighlightColorIndex = $wdYellow
}

# 2. "All names are unique" bullet (the one describing Gym names, which
#    follows the "Gym name cannot be null or empty." validation bullet) ->
#    highlight the whole paragraph, including the paragraph mark.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq "All names are unique`r") {
        $prev = $paras.Item($i - 1)
        if ($prev.Range.Text -like "*Gym name cannot be null or empty*") {
            $p.Range.Font.HighlightColorIndex = $wdYellow
        }
    }
}

# 3. "The Gym name passed to the methods will always be valid!" -> highlight
#    everything from "The" through "valid!" but not the space preceding it.
$range2 = $d.Content
$found2 = $range2.Find.Execute("The Gym name passed to the methods will always be valid!")
if ($found2) {
    $range2.Font.HighlightColorIndex = $wdYellow
}
